# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 19:52"

# --- Update case numbers for several countries (columns B:H = Casos totales,
#     Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 866148
$ws.Range("C4").Value = 17431
$ws.Range("E4").Value = 732468
$ws.Range("G4").Value = 1209
$ws.Range("H4").Value = 48868

# Row 8: Alemania
$ws.Range("B8").Value = 151285
$ws.Range("C8").Value = 637
$ws.Range("E8").Value = 42618
$ws.Range("G8").Value = 52
$ws.Range("H8").Value = 5367

# Row 16: Canada
$ws.Range("B16").Value = 41791
$ws.Range("C16").Value = 1601
$ws.Range("E16").Value = 25777

# Row 18: Suiza
$ws.Range("E18").Value = 7047
$ws.Range("G18").Value = 40
$ws.Range("H18").Value = 1549

# Row 31: Pakistan
$ws.Range("B31").Value = 10880
$ws.Range("C31").Value = 804
$ws.Range("E31").Value = 8315

# Row 108: Jordania
$ws.Range("B108").Value = 437
$ws.Range("C108").Value = 2
$ws.Range("D108").Value = 318
$ws.Range("E108").Value = 112

# Row 127: Paraguay
$ws.Range("D127").Value = 67
$ws.Range("E127").Value = 137

# --- Reorder "Sierra Leona" ahead of "Uganda"/"Haiti" (rows 158-160), with
#     Sierra Leona's numbers updated and Uganda/Haiti shifted down one row ---
$ws.Range("A158").Value = "Sierra Leona"
$ws.Range("B158").Value = 64
$ws.Range("C158").Value = 3
$ws.Range("D158").Value = 10
$ws.Range("E158").Value = 53
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 1

$ws.Range("A159").Value = "Uganda"
$ws.Range("B159").Value = 63
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 46
$ws.Range("E159").Value = 17
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0

$ws.Range("A160").Value = "Haiti"
$ws.Range("B160").Value = 62
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 2
$ws.Range("E160").Value = 56
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 4

# --- Reorder "Nicaragua" ahead of "Groenlandia"/"Gambia"/"Surinam" (rows 200-203),
#     with Nicaragua's numbers updated and the others shifted down one row ---
$ws.Range("A200").Value = "Nicaragua"
$ws.Range("B200").Value = 11
$ws.Range("C200").Value = 1
$ws.Range("D200").Value = 7
$ws.Range("E200").Value = 1
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 1
$ws.Range("H200").Value = 3

$ws.Range("A201").Value = "Groenlandia"
$ws.Range("B201").Value = 11
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 11
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

$ws.Range("A202").Value = "Gambia"
$ws.Range("B202").Value = 10
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 2
$ws.Range("E202").Value = 7
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 1

$ws.Range("A203").Value = "Surinam"
$ws.Range("B203").Value = 10
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 6
$ws.Range("E203").Value = 3
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 1
